$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete trailing rows 15-18 (sheet shrinks from A1:F18 to A1:F14)
$ws.Rows("15:18").Delete()

# Update stock-symbol grid B2:F14 to new values
$ws.Range("B2").Value = "NSE:COROMANDEL"
$ws.Range("C2").Value = "NSE:AARTIDRUGS"
$ws.Range("D2").Value = "NSE:GAIL"
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = "NSE:COROMANDEL"
$ws.Range("B3").Value = "NSE:DHRUV"
$ws.Range("C3").Value = "NSE:ASIANENE"
$ws.Range("D3").Value = "NSE:GODREJPROP"
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = "NSE:DLF"
$ws.Range("B4").Value = "NSE:DLF"
$ws.Range("C4").Value = "NSE:BAJAJELEC"
$ws.Range("D4").Value = "NSE:HDFCAMC"
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = "NSE:HAVELLS"
$ws.Range("B5").Value = "NSE:DODLA"
$ws.Range("C5").Value = "NSE:DCMSHRIRAM"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = "NSE:LAURUSLABS"
$ws.Range("B6").Value = "NSE:HDFCNIFIT"
$ws.Range("C6").Value = "NSE:DOLATALGO"
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = "NSE:LTTS"
$ws.Range("B7").Value = "NSE:ICDSLTD"
$ws.Range("C7").Value = "NSE:GFLLIMITED"
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = "NSE:POWERGRID"
$ws.Range("B8").Value = "NSE:KIMS"
$ws.Range("C8").Value = "NSE:HGS"
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = "NSE:RAMCOCEM"
$ws.Range("B9").Value = "NSE:KRSNAA"
$ws.Range("C9").Value = "NSE:IRMENERGY"
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""
$ws.Range("F9").Value = ""
$ws.Range("B10").Value = "NSE:MASPTOP50"
$ws.Range("C10").Value = "NSE:KAPSTON"
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = ""
$ws.Range("F10").Value = ""
$ws.Range("B11").Value = "NSE:MONQ50"
$ws.Range("C11").Value = "NSE:LAMBODHARA"
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = ""
$ws.Range("B12").Value = "NSE:PAYTM"
$ws.Range("C12").Value = "NSE:MACPOWER"
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = ""
$ws.Range("B13").Value = "NSE:PHOENIXLTD"
$ws.Range("C13").Value = "NSE:POKARNA"
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = ""
$ws.Range("F13").Value = ""
$ws.Range("B14").Value = "NSE:RUSTOMJEE"
$ws.Range("C14").Value = ""
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = ""
$ws.Range("F14").Value = ""
